$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the dates on rows 8 and 9 (44626 -> 44627), keep the existing date style
$ws.Range("A8").Value = 44627
$ws.Range("A9").Value = 44627

# Add the new challenge of the week on rows 10 and 11, copying row 9's
# formatting so the new date cells keep the same number format/style
$ws.Range("A9").Copy()
$ws.Range("A10:A11").PasteSpecial(-4122)
$ws.Range("B9").Copy()
$ws.Range("B10:B11").PasteSpecial(-4122)

$ws.Range("A10").Value = 44634
$ws.Range("B10").Value = "Ecriture aide-mémoire Javascript"

$ws.Range("A11").Value = 44641
$ws.Range("B11").Value = "Ecriture aide-mémoire Javascript"

# Match the saved selection state
$ws.Range("B10:B11").Select()
